$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look numeric need to be forced to Text format
# first, otherwise Excel auto-converts them to floating point numbers and
# the exact original text (trailing zeros, etc.) is lost.
$textCells = @("D5", "D9", "D10", "D17", "D21", "D22", "D26", "D28", "D29", "D31", "D32", "D38", "D39", "D41", "D47", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values
$ws.Range("D2").Value = '25.916.05'
$ws.Range("E2").Value = '  +0.76%  '
$ws.Range("E3").Value = '  +0.14%  '
$ws.Range("E4").Value = '  -0.26%  '
$ws.Range("D5").Value = '214.87'
$ws.Range("E6").Value = '  -0.25%  '
$ws.Range("E8").Value = '  -0.59%  '
$ws.Range("D9").Value = '0.0632'
$ws.Range("E9").Value = '  -0.46%  '
$ws.Range("D10").Value = '19.67'
$ws.Range("E10").Value = '  +0.10%  '
$ws.Range("E11").Value = '  +0.52%  '
$ws.Range("E12").Value = '  +0.59%  '
$ws.Range("B13").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C13").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D13").Value = '1.861.37'
$ws.Range("E13").Value = '  +0.22%  '
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '1.641.97'
$ws.Range("E14").Value = '  +0.95%  '
$ws.Range("E15").Value = '  -0.79%  '
$ws.Range("E16").Value = '  -0.45%  '
$ws.Range("D17").Value = '62.95'
$ws.Range("E17").Value = '  +0.50%  '
$ws.Range("D18").Value = '25.905.98'
$ws.Range("E18").Value = '  +0.64%  '
$ws.Range("E19").Value = '  -0.25%  '
$ws.Range("E20").Value = '  +0.18%  '
$ws.Range("D21").Value = '191.78'
$ws.Range("E21").Value = '  -1.08%  '
$ws.Range("D22").Value = '9.99'
$ws.Range("E23").Value = '  +0.91%  '
$ws.Range("E24").Value = '  -0.23%  '
$ws.Range("E25").Value = '  -1.59%  '
$ws.Range("D26").Value = '142.27'
$ws.Range("E26").Value = '  +0.13%  '
$ws.Range("E27").Value = '  +1.22%  '
$ws.Range("D28").Value = '6.87'
$ws.Range("E28").Value = '  -0.06%  '
$ws.Range("D29").Value = '15.54'
$ws.Range("E29").Value = '  +0.15%  '
$ws.Range("E30").Value = '  +0.06%  '
$ws.Range("D31").Value = '0.0493'
$ws.Range("E31").Value = '  +0.34%  '
$ws.Range("D32").Value = '3.33'
$ws.Range("E32").Value = '  +0.43%  '
$ws.Range("E33").Value = '  +0.23%  '
$ws.Range("E34").Value = '  +1.03%  '
$ws.Range("E35").Value = '  +0.38%  '
$ws.Range("E36").Value = '  +0.89%  '
$ws.Range("D37").Value = '1.147.36'
$ws.Range("E37").Value = '  +1.89%  '
$ws.Range("D38").Value = '0.545'
$ws.Range("E38").Value = '  -0.22%  '
$ws.Range("D39").Value = '2.51'
$ws.Range("E39").Value = '  -0.79%  '
$ws.Range("E40").Value = '  +0.76%  '
$ws.Range("D41").Value = '1.00'
$ws.Range("E41").Value = '  -0.26%  '
$ws.Range("E42").Value = '  +1.40%  '
$ws.Range("E43").Value = '  +1.16%  '
$ws.Range("E44").Value = '  -0.02%  '
$ws.Range("D45").Value = '1.771.99'
$ws.Range("E45").Value = '  +0.27%  '
$ws.Range("D46").Value = '0.0₆0110'
$ws.Range("E46").Value = '  -1.74%  '
$ws.Range("D47").Value = '55.61'
$ws.Range("E47").Value = '  +1.11%  '
$ws.Range("E48").Value = '  +1.86%  '
$ws.Range("E49").Value = '  +5.66%  '
$ws.Range("E50").Value = '  +0.03%  '
$ws.Range("D51").Value = '7.59'
$ws.Range("E51").Value = '  +0.83%  '
